$wb = $excel.ActiveWorkbook

# --- Section_A sheet updates ---
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("C2").Value = "Free"
$wsA.Range("D2").Value = "DS456"
$wsA.Range("F2").Value = "DS456"

$wsA.Range("C3").Value = "Free"
$wsA.Range("D3").Value = "Free"

$wsA.Range("B5").Value = "Free"

$wsA.Range("B6").Value = "DS401"
$wsA.Range("C6").Value = "DS456"
$wsA.Range("E6").Value = "Free"
$wsA.Range("F6").Value = "DS401"

$wsA.Range("C7").Value = "EC456"
$wsA.Range("D7").Value = "EC456"
$wsA.Range("E7").Value = "DS456"
$wsA.Range("F7").Value = "EC456"

# --- Section_B sheet updates ---
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "DS401"
$wsB.Range("C2").Value = "Free"
$wsB.Range("D2").Value = "Free"
$wsB.Range("E2").Value = "DS401"

$wsB.Range("B3").Value = "EC456"
$wsB.Range("C3").Value = "Free"
$wsB.Range("D3").Value = "EC456"
$wsB.Range("E3").Value = "Free"
$wsB.Range("F3").Value = "DS456"

$wsB.Range("B5").Value = "DS456"
$wsB.Range("C5").Value = "DS456"
$wsB.Range("D5").Value = "DS456"
$wsB.Range("E5").Value = "Free"
$wsB.Range("F5").Value = "DS401"

$wsB.Range("B6").Value = "Free"
$wsB.Range("D6").Value = "DS401"

$wsB.Range("C7").Value = "EC456"
